$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '71.414.07'
$ws.Range("E2").Value = '  +0.94%  '

$ws.Range("D3").Value = '3.642.23'
$ws.Range("E3").Value = '  +0.31%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  +0.11%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '582.72'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -1.16%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '189.66'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.04%  '

$ws.Range("D7").Value = '3.639.72'
$ws.Range("E7").Value = '  +0.38%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.628'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -1.92%  '

$ws.Range("E9").Value = '  +0.09%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.189'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +3.48%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.662'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.31%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '55.19'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -4.68%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000309'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.27%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '9.76'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -1.25%  '

$ws.Range("D15").Value = '4.237.08'
$ws.Range("E15").Value = '  +0.54%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '19.83'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -2.27%  '

$ws.Range("D17").Value = '3.653.47'
$ws.Range("E17").Value = '  +0.56%  '

$ws.Range("D18").Value = '71.397.13'
$ws.Range("E18").Value = '  +0.94%  '

$ws.Range("E19").Value = '  -0.10%  '

$ws.Range("E20").Value = '  -0.91%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '1.06'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.66%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '506.48'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.93%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '19.51'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +1.50%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '4.96'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -3.23%  '

$ws.Range("E25").Value = '  -0.78%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '96.64'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +6.58%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '11.69'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.47%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '3.03'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -3.82%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '9.48'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -0.26%  '

$ws.Range("E30").Value = '  -2.69%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.23'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -1.45%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '12.74'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.85%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '66.37'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -0.84%  '

$ws.Range("E34").Value = '  -2.28%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '579.70'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -5.25%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '3.26'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +9.51%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '39.40'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -2.45%  '

$ws.Range("E38").Value = '  +1.52%  '

$ws.Range("D39").Value = '0.0₃0808'
$ws.Range("E39").Value = '  -3.54%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.999'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.03%  '

$ws.Range("E41").Value = '  +5.44%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.51'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -2.42%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.138'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -7.31%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '3.09'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -2.24%  '

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0458'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.22%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.59'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.77%  '

$ws.Range("D47").Value = '3.246.77'
$ws.Range("E47").Value = '  -2.41%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '9.68'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.42%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.137'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -1.11%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '1.54'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +28.14%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.998'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -0.15%  '
